$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) cells are stored as text, matching the source data
# (many price values look numeric, e.g. "1.00", and must not be auto-converted)
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '60.139.65'
$ws.Range("E2").Value = '  -3.03%  '

$ws.Range("D3").Value = '3.290.03'
$ws.Range("E3").Value = '  -3.89%  '

$ws.Range("D4").Value = '1.00'
$ws.Range("E4").Value = '  +0.02%  '

$ws.Range("D5").Value = '556.19'
$ws.Range("E5").Value = '  -3.93%  '

$ws.Range("D6").Value = '141.25'
$ws.Range("E6").Value = '  -8.01%  '

$ws.Range("E7").Value = '  -0.10%  '

$ws.Range("D8").Value = '3.291.62'
$ws.Range("E8").Value = '  -3.83%  '

$ws.Range("D9").Value = '0.466'
$ws.Range("E9").Value = '  -3.92%  '

$ws.Range("D10").Value = '7.93'
$ws.Range("E10").Value = '  -1.38%  '

$ws.Range("E11").Value = '  -5.57%  '

$ws.Range("D12").Value = '0.406'
$ws.Range("E12").Value = '  -3.07%  '

$ws.Range("D13").Value = '3.850.72'
$ws.Range("E13").Value = '  -3.97%  '

$ws.Range("E14").Value = '  -0.06%  '

$ws.Range("D15").Value = '26.68'
$ws.Range("E15").Value = '  -6.32%  '

$ws.Range("D16").Value = '3.284.74'
$ws.Range("E16").Value = '  -4.38%  '

$ws.Range("E17").Value = '  -5.05%  '

$ws.Range("D18").Value = '60.165.13'
$ws.Range("E18").Value = '  -3.01%  '

$ws.Range("D19").Value = '6.07'
$ws.Range("E19").Value = '  -7.50%  '

$ws.Range("D20").Value = '13.76'
$ws.Range("E20").Value = '  -5.15%  '

$ws.Range("D21").Value = '8.51'
$ws.Range("E21").Value = '  -5.03%  '

$ws.Range("D22").Value = '371.94'
$ws.Range("E22").Value = '  -2.90%  '

$ws.Range("E23").Value = '  +0.07%  '

$ws.Range("D24").Value = '72.29'
$ws.Range("E24").Value = '  -4.16%  '

$ws.Range("D25").Value = '0.532'
$ws.Range("E25").Value = '  -7.06%  '

$ws.Range("D26").Value = '3.413.30'
$ws.Range("E26").Value = '  -4.13%  '

$ws.Range("E27").Value = '  -9.12%  '

$ws.Range("D28").Value = '0.172'
$ws.Range("E28").Value = '  -3.57%  '

$ws.Range("E29").Value = '  -0.15%  '

$ws.Range("D30").Value = '7.06'
$ws.Range("E30").Value = '  -7.86%  '

$ws.Range("E31").Value = '  -0.01%  '

$ws.Range("E32").Value = '  -5.44%  '

$ws.Range("E33").Value = '  -6.00%  '

$ws.Range("D34").Value = '22.52'
$ws.Range("E34").Value = '  -3.25%  '

$ws.Range("E35").Value = '  -7.47%  '

$ws.Range("B36").Value = 'NEARProtocol'
$ws.Range("C36").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D36").Value = '5.06'
$ws.Range("E36").Value = '  -8.81%  '

$ws.Range("B37").Value = 'Monero'
$ws.Range("C37").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range("D37").Value = '166.10'
$ws.Range("E37").Value = '  -1.43%  '

$ws.Range("E38").Value = '  -5.51%  '

$ws.Range("D39").Value = '6.61'
$ws.Range("E39").Value = '  -5.23%  '

$ws.Range("D40").Value = '3.319.47'
$ws.Range("E40").Value = '  -3.98%  '

$ws.Range("E41").Value = '  -8.19%  '

$ws.Range("D42").Value = '25.79'
$ws.Range("E42").Value = '  -17.18%  '

$ws.Range("D43").Value = '41.53'
$ws.Range("E43").Value = '  -2.62%  '

$ws.Range("D44").Value = '0.745'
$ws.Range("E44").Value = '  -4.58%  '

$ws.Range("B45").Value = 'ONDO'
$ws.Range("C45").Value = 'https://coinranking.com/coin/7AQlxzQpQ+ondo-ondo'
$ws.Range("D45").Value = '1.12'
$ws.Range("E45").Value = '  -3.78%  '

$ws.Range("B46").Value = 'Filecoin'
$ws.Range("C46").Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range("D46").Value = '4.10'
$ws.Range("E46").Value = '  -7.78%  '

$ws.Range("E47").Value = '  -7.08%  '

$ws.Range("E48").Value = '  -0.03%  '

$ws.Range("D49").Value = '2.315.19'
$ws.Range("E49").Value = '  -9.33%  '

$ws.Range("D50").Value = '6.33'
$ws.Range("E50").Value = '  -7.26%  '

$ws.Range("D51").Value = '21.43'
$ws.Range("E51").Value = '  -5.80%  '
